# Edit script: add "Country Name" column, populate US (existing rows) and
# Japan (two new rows) for the World Development Indicators extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column A. This shifts the existing
#    data (old A:AV) one column to the right (new B:AW), preserving all
#    existing shared-string references, styles and values untouched.
$ws.Columns("A").Insert()

# 2. Populate the new column A: header + "US" for the pre-existing rows.
$ws.Range("A1").Value2 = "Country Name"
$ws.Range("A2").Value2 = "US"
$ws.Range("A3").Value2 = "US"
$ws.Range("A4").Value2 = "US"
$ws.Range("A5").Value2 = "US"
$ws.Range("A6").Value2 = "US"

# 3. Add two new rows (7 and 8) of data for Japan, re-using the same two
#    series already present for the US (rows 2 and 3: GDP per capita and
#    NPISH final consumption expenditure).
$ws.Range("A7").Value2 = "Japan"
$ws.Range("B7").Value2 = $ws.Range("B2").Value2
$ws.Range("A8").Value2 = "Japan"
$ws.Range("B8").Value2 = $ws.Range("B3").Value2

$row7 = @(9456.634807052107,9627.71231827424,9920.345697738203,10346.65830848554,10925.208306801964,10955.1528376292,11149.901957858696,11588.994335694033,11897.427194644573,12192.809053272498,12613.528687707212,13007.846876641212,13511.650000918022,14147.911733337116,14785.369005817083,15438.911186799216,15719.740278312149,16017.824540144065,16137.978294896026,16460.631774959715,16826.417235556648,17117.04205921647,17170.987608860287,17013.6273260945,17163.920266750312,17399.11291824812,17713.177192737163,17899.76924226264,17974.771779191924,18206.19184554176,18483.12907259045,18642.9046524387,18762.11899676706,18540.530274656234,18373.162174374862,18797.965609383395,18734.040857077664,19144.624779475762,19671.919648031577,19519.46069445279,19497.82710219738,19424.56555403069,19645.066765817868,19715.403242230488,19619.976630659556,18803.014638046352,19037.611669264974)
$row8 = @(14102.69164648237,14647.06439506061,15214.42996771421,16056.819488225376,16389.53593245061,16290.424149130547,16700.08868366794,17096.353482805956,17437.215881929984,18154.871243954967,19166.238861775044,19884.439703440385,20570.858676414464,21998.699775847726,22887.451473595123,23891.10267896754,24367.943527404284,24237.680154175097,23942.976067639513,24433.95568322401,25142.446590664673,25945.625810720278,25903.69151674802,25355.317814312202,25082.844868648728,25726.833979553525,25795.424979482243,25619.500713500132,25978.033586500555,26597.914256688135,26934.849279392954,27051.69252071199,27312.114104509976,26020.688158929937,24510.10748245739,25719.781371293426,25572.95362083791,26015.128274636525,26755.657775908057,26886.020643449483,28006.306867406423,28412.549489899124,28832.862637905753,28762.416500272513,28694.791597665997,27232.521434157494,27008.840862738867)

for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 3).Value2 = $row7[$i]
    $ws.Cells.Item(8, $i + 3).Value2 = $row8[$i]
}

# Match the number format used by the other value rows (style index 3,
# format code "0") for the two new data rows.
$ws.Range("C7:AW8").NumberFormat = "0"

# 4. Column widths: the new column A is narrower (16 chars) than the old
#    series-name column (now B), and the new last column (AW, holding the
#    wide "Country Name"/"US"/"Japan" header text) needs to widen too.
$ws.Columns("A").ColumnWidth = 15.17

# 5. Sheet view: selection moves to B12, and the saved "topLeftCell" scroll
#    position is cleared (back to the sheet's natural A1 origin).
$ws.Range("B12").Select()
